$d = $word.ActiveDocument

$replacements = @(
    @{old="954÷7=136, 2"; new="216÷8=27, 0"},
    @{old="946÷8=118, 2"; new="740÷8=92, 4"},
    @{old="703÷9=78, 1"; new="681÷8=85, 1"},
    @{old="526÷5=105, 1"; new="484÷9=53, 7"},
    @{old="517÷2=258, 1"; new="858÷3=286, 0"},
    @{old="648÷3=216, 0"; new="342÷8=42, 6"},
    @{old="307÷7=43, 6"; new="620÷3=206, 2"},
    @{old="512÷2=256, 0"; new="568÷4=142, 0"},
    @{old="275÷9=30, 5"; new="780÷8=97, 4"},
    @{old="438÷3=146, 0"; new="581÷8=72, 5"},
    @{old="529÷3=176, 1"; new="469÷2=234, 1"},
    @{old="280÷6=46, 4"; new="897÷8=112, 1"},
    @{old="444÷6=74, 0"; new="104÷6=17, 2"},
    @{old="584÷4=146, 0"; new="674÷6=112, 2"},
    @{old="719÷8=89, 7"; new="758÷3=252, 2"},
    @{old="690÷4=172, 2"; new="259÷2=129, 1"},
    @{old="636÷2=318, 0"; new="809÷6=134, 5"},
    @{old="113÷8=14, 1"; new="373÷2=186, 1"},
    @{old="368÷8=46, 0"; new="486÷3=162, 0"},
    @{old="631÷5=126, 1"; new="823÷2=411, 1"},
    @{old="433÷6=72, 1"; new="695÷8=86, 7"},
    @{old="563÷6=93, 5"; new="575÷7=82, 1"},
    @{old="430÷6=71, 4"; new="677÷8=84, 5"},
    @{old="257÷9=28, 5"; new="366÷7=52, 2"},
    @{old="638÷9=70, 8"; new="898÷6=149, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
